# RandomFloat_InsertSortTimes.csv.xlsx - "Fix sorting and generate viable
# xlsx and charts":
#
# The average timings for the two smallest sample sizes were recomputed,
# changing the Avg_Time_ms values (column D) for the 5000- and 10000-row
# runs:
#   D2 (Rows=5000) : 17.847199 -> 18.229861
#   D3 (Rows=10000): 72.439356 -> 77.784881
#
# The worksheet's scatter chart plots Avg_Time_ms (Data!$D$2:$D$8) as its
# y-values, so updating these two cells is the source-of-truth edit; the
# chart is refreshed afterwards to keep it in sync with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 18.229861
$ws.Range("D3").Value = 77.784881

try {
    $ws.ChartObjects().Item(1).Chart.Refresh()
} catch {
}
